$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3555633.88
$ws.Range("C9").Value = 557560.24
$ws.Range("D9").Value = 4113194.12
$ws.Range("E9").Value = 13.55540788335076
$ws.Range("F9").Value = 86.44459211664923
$ws.Range("G9").Value = -46.11458292752495
$ws.Range("H9").Value = -35.79016686520371
$ws.Range("I9").Value = 35761
$ws.Range("J9").Value = 1524
$ws.Range("K9").Value = 37285
$ws.Range("L9").Value = 25726
$ws.Range("M9").Value = 159.8847127419731
$ws.Range("N9").Value = 9.156549877510933

$wb.Save()
